$wb = $excel.ActiveWorkbook

# --- Sheet "Schedule" updates ---
$wsSchedule = $wb.Worksheets.Item("Schedule")
$wsSchedule.Range("E3").Value = 536.0526990000001
$wsSchedule.Range("F3").Value = 35.45322083333334
$wsSchedule.Range("E4").Value = -99.12418425000001
$wsSchedule.Range("F4").Value = -2.185277430555556

# --- Sheet "Detailed" updates ---
$wsDetailed = $wb.Worksheets.Item("Detailed")
$wsDetailed.Range("B39").Value = -3.03118
$wsDetailed.Range("B40").Value = -3.04997
$wsDetailed.Range("B41").Value = 75.07939
$wsDetailed.Range("C41").Value = "historical"
$wsDetailed.Range("B42").Value = 112.59605
$wsDetailed.Range("C42").Value = "historical"
$wsDetailed.Range("B43").Value = 77.18000000000001
$wsDetailed.Range("B44").Value = 66.99298
$wsDetailed.Range("B45").Value = 57.04922
$wsDetailed.Range("B46").Value = 64.8901
$wsDetailed.Range("B47").Value = 57.04922
$wsDetailed.Range("B48").Value = 57.06007
$wsDetailed.Range("B49").Value = 56.98
$wsDetailed.Range("B51").Value = 42.92029
$wsDetailed.Range("B52").Value = 36.06
$wsDetailed.Range("B54").Value = 36.06
$wsDetailed.Range("B57").Value = 35.88
$wsDetailed.Range("B58").Value = 35.88
$wsDetailed.Range("B59").Value = 56.98
$wsDetailed.Range("B60").Value = 53.15912
$wsDetailed.Range("B61").Value = 56.98
$wsDetailed.Range("B64").Value = 36.06
$wsDetailed.Range("B65").Value = 0.7
$wsDetailed.Range("B66").Value = -4.80722
$wsDetailed.Range("B67").Value = -4.81333
$wsDetailed.Range("B68").Value = 0.7
$wsDetailed.Range("B69").Value = 0.59034
$wsDetailed.Range("B70").Value = 0.009719999999999999
$wsDetailed.Range("B71").Value = 0.00003
$wsDetailed.Range("B72").Value = 0.00002
$wsDetailed.Range("B73").Value = 0
$wsDetailed.Range("B74").Value = -5.01
$wsDetailed.Range("B75").Value = -5.39446
$wsDetailed.Range("B76").Value = -1.74235
$wsDetailed.Range("B77").Value = -6.8
$wsDetailed.Range("B78").Value = -14
$wsDetailed.Range("B80").Value = -14
$wsDetailed.Range("B81").Value = -5.68475
$wsDetailed.Range("B82").Value = 5.27701
$wsDetailed.Range("B83").Value = -6.41818
$wsDetailed.Range("B84").Value = -7.67416
$wsDetailed.Range("B85").Value = -6.39013
$wsDetailed.Range("B86").Value = -6.20642
$wsDetailed.Range("B87").Value = -6.01284
$wsDetailed.Range("B88").Value = 0.01089
$wsDetailed.Range("B89").Value = 2.23907
$wsDetailed.Range("B90").Value = 17.98199
$wsDetailed.Range("B91").Value = 8.49396
$wsDetailed.Range("B92").Value = 6.10154
$wsDetailed.Range("B93").Value = 4.12106
$wsDetailed.Range("B94").Value = 36.05728
$wsDetailed.Range("B95").Value = 57.3
$wsDetailed.Range("B96").Value = 57.3
$wsDetailed.Range("B97").Value = 57.06005